# five-sound-fortune (Version 1).docx edit script
#
# What this does (per the target diff):
#   1. Removes the "Meta description: Read our review of Five Sound Fortune
#      online slot game and play for free. Learn about the gameplay,
#      symbols, and more." paragraph that originally sat right under the
#      title.
#   2. Inserts a new bold paragraph, "Play Five Sound Fortune for Free -
#      Review of the Online Slot Game", right before the final (italic)
#      paragraph at the end of the document.
#   3. Replaces the text of that final italic paragraph - which used to be
#      the "Create a feature image..." image prompt - with "Read our review
#      of Five Sound Fortune online slot game and play for free. Learn
#      about the gameplay, symbols, and more." (keeping its italic run
#      formatting intact).

$d = $word.ActiveDocument

# --- Step 1: delete the "Meta description" paragraph (2nd paragraph). -----
$d.Paragraphs.Item(2).Range.Delete() | Out-Null

# --- Step 2: insert a new bold paragraph just before the last paragraph. --
$count = $d.Paragraphs.Count
$lastRange = $d.Paragraphs.Item($count).Range
$insertionPoint = $d.Range($lastRange.Start, $lastRange.Start)

# Build a FlatOPC OOXML fragment containing the new paragraph. A trailing
# empty <w:p/> is included so the paragraph break lands cleanly before the
# existing final paragraph instead of merging into it; the stray empty
# paragraph that results is removed right afterwards.
$newParaXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Five Sound Fortune for Free - Review of the Online Slot Game</w:t></w:r></w:p><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$insertionPoint.InsertXML($newParaXml) | Out-Null

# Remove the placeholder empty paragraph InsertXML left behind.
$count2 = $d.Paragraphs.Count
$d.Paragraphs.Item($count2 - 1).Range.Delete() | Out-Null

# --- Step 3: swap the final paragraph's text for the new description. -----
$oldText = "Create a feature image that captures the fun and festive atmosphere of " + [char]34 + "Five Sound Fortune" + [char]34 + ". The image should be in cartoon style and feature a happy Maya warrior with glasses. The Maya warrior should be holding a red lantern and surrounded by fireworks. The background should depict a street with illuminated houses, as described in the game review. The colors should be bright and vibrant to reflect the celebratory theme of the game. Use your creativity to bring the game to life and entice players to give it a spin!"
$newText = "Read our review of Five Sound Fortune online slot game and play for free. Learn about the gameplay, symbols, and more."

$d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null
